$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from B1 to C1, then set new header text
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = 'Coord: normal vector scan'

# Update column B values (recomputed) and add column C vector strings
$ws.Range("B2").Value = 0.1333730812976081
$ws.Range("C2").Value = '[0.         0.32276389 0.94647952]'
$ws.Range("B3").Value = 1.611874410507632
$ws.Range("C3").Value = '[-0.43690781  0.52136285  0.73300228]'
$ws.Range("B4").Value = 0.3267080566507467
$ws.Range("C4").Value = '[-0.00531186  0.01811423  0.99982181]'
$ws.Range("B5").Value = 0.7858762146602799
$ws.Range("C5").Value = '[-6.91565831e-04  2.81342319e-01 -9.59607222e-01]'
$ws.Range("B6").Value = 2.087626368550382
$ws.Range("C6").Value = '[0.71928458 0.28914187 0.63168557]'
$ws.Range("B7").Value = 0.2261021724828858
$ws.Range("C7").Value = '[-0.73973604 -0.27273953  0.6151453 ]'
$ws.Range("B8").Value = 0.4609665665156076
$ws.Range("C8").Value = '[0.         0.31734706 0.94830946]'
$ws.Range("B9").Value = 1.499799544499535
$ws.Range("C9").Value = '[ 0.         -0.30010201  0.95390712]'
$ws.Range("B10").Value = 2.038754803392189
$ws.Range("C10").Value = '[-0.72578655  0.25637797  0.63836057]'
$ws.Range("B11").Value = 2.211307404105731
$ws.Range("C11").Value = '[ 0.72576295 -0.27321059  0.63136686]'
$ws.Range("B12").Value = 0.7699691685337028
$ws.Range("C12").Value = '[ 6.77333416e-04 -2.81608377e-01 -9.59529188e-01]'
$ws.Range("B13").Value = 2.393526889858017
$ws.Range("C13").Value = '[-0.7179608  -0.28315001  0.63589179]'
$ws.Range("B14").Value = 4.414422470166294
$ws.Range("C14").Value = '[0.69516985 0.25705402 0.67131372]'
$ws.Range("B15").Value = 1.614874150419871
$ws.Range("C15").Value = '[ 0.         -0.29818554  0.95450793]'
